$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Fix typo in the "Digerença" -> "Diferença" header before moving things around
$ws1.Range("N2").Value = "Diferença"

# Create Sheet2 positioned after Sheet1 (becomes the new active sheet)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Move the "Tabela2" block (J2:N8 on Sheet1) onto Sheet2 at A1:E7
$ws1.Range("J2:N8").Cut($ws2.Range("A1"))

# Drop the old table definition bound to Sheet1 and rebuild it on Sheet2
$lo = $ws1.ListObjects.Item("Tabela2")
$lo.Unlist()

$lo2 = $ws2.ListObjects.Add(1, $ws2.Range("A1:E6"), [System.Reflection.Missing]::Value, 1)
$lo2.Name = "Tabela2"
$lo2.TableStyle = "Sheet1-style"

# Restore the formulas that were flattened to values by the cross-sheet Cut
$ws2.Range("B2").Formula = "=AVERAGE(Table1[test_accuracy])"
$ws2.Range("C2").Formula = "=STDEV(Table1[test_accuracy])"
$ws2.Range("E2").Value = "n/a"

$ws2.Range("B3").Formula = "=AVERAGE(Table1[test_recall])"
$ws2.Range("C3").Formula = "=STDEV(Table1[test_recall])"
$ws2.Range("E3").Formula = "=(Tabela2[[#This Row],[Kvasir-SEG]]-Tabela2[[#This Row],[Média]])/Tabela2[[#This Row],[Kvasir-SEG]]"

$ws2.Range("B4").Formula = "=AVERAGE(Table1[test_precision])"
$ws2.Range("C4").Formula = "=STDEV(Table1[test_precision])"
$ws2.Range("E4").Formula = "=Tabela2[[#This Row],[Kvasir-SEG]]-Tabela2[[#This Row],[Média]]"

$ws2.Range("B5").Formula = "=AVERAGE(Table1[test_iou])"
$ws2.Range("C5").Formula = "=STDEV(Table1[test_iou])"
$ws2.Range("E5").Formula = "=Tabela2[[#This Row],[Kvasir-SEG]]-Tabela2[[#This Row],[Média]]"

$ws2.Range("B6").Formula = "=AVERAGE(Table1[test_f1])"
$ws2.Range("C6").Formula = "=STDEV(Table1[test_f1])"
$ws2.Range("E6").Formula = "=Tabela2[[#This Row],[Kvasir-SEG]]-Tabela2[[#This Row],[Média]]"

# Totals row
$lo2.ShowTotals = $true
$lo2.ListColumns.Item(5).TotalsCalculation = 9
$ws2.Range("E7").Formula = "=MEDIAN(E3:E6)"

Write-Host "done"
